$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.569988608360291
$ws.Range("B1").Value = 2.426857233047485
$ws.Range("C1").Value = 5.450758457183838
$ws.Range("D1").Value = 1.493471145629883
$ws.Range("E1").Value = 0.8236287236213684
